$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productos")

# --- Update STOCK (E) / PRECIO... actually columns are: A=ID, B=DESCRIPCION, C=STOCK, D=PRECIO, E=CATEGORIA, F=PROVEEDOR
# Column E/F values change per diff for rows 2-20, plus D16/D17 swap, and a new row 21.

# Rows 2-6: E 5->3, F 3->1
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 5).Value = 3
    $ws.Cells.Item($r, 6).Value = 1
}

# Rows 7-12: E 3->1, F 1->2
for ($r = 7; $r -le 12; $r++) {
    $ws.Cells.Item($r, 5).Value = 1
    $ws.Cells.Item($r, 6).Value = 2
}

# Rows 13-15: E 3->1, F 1->3
for ($r = 13; $r -le 15; $r++) {
    $ws.Cells.Item($r, 5).Value = 1
    $ws.Cells.Item($r, 6).Value = 3
}

# Row 16: D 100->150, E 4->2, F 3->4
$ws.Cells.Item(16, 4).Value = 150
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 4

# Row 17: D 150->100, E 4->2, F 3->4
$ws.Cells.Item(17, 4).Value = 100
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 4

# Row 18: E 6->4, F 3->5
$ws.Cells.Item(18, 5).Value = 4
$ws.Cells.Item(18, 6).Value = 5

# Rows 19-20: E 6->4, F 3->6
$ws.Cells.Item(19, 5).Value = 4
$ws.Cells.Item(19, 6).Value = 6
$ws.Cells.Item(20, 5).Value = 4
$ws.Cells.Item(20, 6).Value = 6

# New row 21: Cadbury Tableta x80g
$ws.Cells.Item(21, 1).Value = 59
$ws.Cells.Item(21, 2).Value = "Cadbury Tableta x80g"
$ws.Cells.Item(21, 3).Value = 300
$ws.Cells.Item(21, 4).Value = 125
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 6).Value = 5

# Column B width (widened to fit the longer product descriptions)
$ws.Columns.Item(2).ColumnWidth = 27.31

# Selection changes to B20
$ws.Range("B20").Select()

$wb.Save()
